$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three comments in column C whose text actually changed.
# (Other cells in column C keep the same text; Excel will renumber the
# shared-string table automatically as strings are added/removed.)

$ws.Range("C10").Value2 = "klart för de ripinventerade lyorna, klart för närmsta vattenkälla för alla lyor. Måste ta avstånd till större vatten. Bestäm ett gränsvärde (m^2), lägg till i attributes och ta bort."

$ws.Range("C11").Value2 = "Fått skript av Rasmus, ändra om till mina data."

$ws.Range("C18").Value2 = "Den högsta sannolikheten för lämmel i en 500 x 500 m pixel är 0,4 under ett uppgångsår. Jag räknar allt från 0,2 och över som bra lämmelhabitat. Toppåren är inte lika intressanta eftersom lämlarna finns överallt."

# Move the active selection from C14 to C20, as it was when the file was
# last saved by the author.
$ws.Range("C20").Select() | Out-Null
